$d = $word.ActiveDocument

$pairs = @(
    @("778÷7=", "935÷2="),
    @("768÷9=", "209÷8="),
    @("216÷4=", "624÷7="),
    @("456÷5=", "638÷4="),
    @("733÷8=", "167÷2="),
    @("976÷4=", "124÷3="),
    @("926÷2=", "516÷3="),
    @("425÷9=", "330÷5="),
    @("952÷5=", "304÷4="),
    @("909÷7=", "173÷6="),
    @("700÷9=", "767÷8="),
    @("775÷6=", "888÷8="),
    @("990÷5=", "111÷2="),
    @("938÷2=", "866÷9="),
    @("281÷5=", "507÷5="),
    @("522÷9=", "881÷8="),
    @("570÷7=", "468÷8="),
    @("805÷6=", "599÷6="),
    @("860÷4=", "843÷3="),
    @("888÷9=", "371÷8="),
    @("572÷9=", "560÷6="),
    @("127÷9=", "232÷5="),
    @("703÷9=", "212÷7="),
    @("719÷4=", "614÷2="),
    @("157÷9=", "905÷7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
